# Generate Report for Handback
# Refresh the handback timestamps on the report for the first file
# (3b55afe9-fca1-4f9d-b841-af26fca2fc20) now that a new handback was
# generated.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: row 2 is the 3b55afe9... file.
$zhcn.Range("H2").Value = "2016-09-06 09:03:36"
$zhcn.Range("K2").Value = "2016-09-06 09:04:48"

# de-de sheet: row 2 is the 3b55afe9... file.
$dede.Range("H2").Value = "2016-09-06 09:03:46"
$dede.Range("K2").Value = "2016-09-06 09:05:15"

# Overview sheet: row 2 is the 3b55afe9... file; "Latest HO Xliff Generate
# Date" reflects the newest of the per-language handoff/handback dates.
$overview.Range("G2").Value = "2016-09-06 09:03:46"
